$d = $word.ActiveDocument

$d.Content.Find.Execute("4+91=", $true, $false, $false, $false, $false, $true, 1, $false, "33-9=", 2) | Out-Null
$d.Content.Find.Execute("69-55=", $true, $false, $false, $false, $false, $true, 1, $false, "79-55=", 2) | Out-Null
$d.Content.Find.Execute("36+53=", $true, $false, $false, $false, $false, $true, 1, $false, "14-5=", 2) | Out-Null
$d.Content.Find.Execute("77-45=", $true, $false, $false, $false, $false, $true, 1, $false, "2+81=", 2) | Out-Null
$d.Content.Find.Execute("94-25=", $true, $false, $false, $false, $false, $true, 1, $false, "36-17=", 2) | Out-Null
$d.Content.Find.Execute("20+64=", $true, $false, $false, $false, $false, $true, 1, $false, "99-6=", 2) | Out-Null
$d.Content.Find.Execute("27+69=", $true, $false, $false, $false, $false, $true, 1, $false, "25+74=", 2) | Out-Null
$d.Content.Find.Execute("42-25=", $true, $false, $false, $false, $false, $true, 1, $false, "27+35=", 2) | Out-Null
$d.Content.Find.Execute("93-19=", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=", 2) | Out-Null
$d.Content.Find.Execute("96-49=", $true, $false, $false, $false, $false, $true, 1, $false, "57-42=", 2) | Out-Null
$d.Content.Find.Execute("44+12=", $true, $false, $false, $false, $false, $true, 1, $false, "4+89=", 2) | Out-Null
$d.Content.Find.Execute("97-85=", $true, $false, $false, $false, $false, $true, 1, $false, "14+5=", 2) | Out-Null
$d.Content.Find.Execute("7+38=", $true, $false, $false, $false, $false, $true, 1, $false, "79+11=", 2) | Out-Null
$d.Content.Find.Execute("18+34=", $true, $false, $false, $false, $false, $true, 1, $false, "27+67=", 2) | Out-Null
$d.Content.Find.Execute("90-4=", $true, $false, $false, $false, $false, $true, 1, $false, "17+17=", 2) | Out-Null
$d.Content.Find.Execute("93-60=", $true, $false, $false, $false, $false, $true, 1, $false, "97-78=", 2) | Out-Null
$d.Content.Find.Execute("31+68=", $true, $false, $false, $false, $false, $true, 1, $false, "97-21=", 2) | Out-Null
$d.Content.Find.Execute("13-2=", $true, $false, $false, $false, $false, $true, 1, $false, "26-24=", 2) | Out-Null
$d.Content.Find.Execute("24+1=", $true, $false, $false, $false, $false, $true, 1, $false, "79-31=", 2) | Out-Null
$d.Content.Find.Execute("4+52=", $true, $false, $false, $false, $false, $true, 1, $false, "84-4=", 2) | Out-Null
$d.Content.Find.Execute("13+23=", $true, $false, $false, $false, $false, $true, 1, $false, "61+9=", 2) | Out-Null
$d.Content.Find.Execute("91-74=", $true, $false, $false, $false, $false, $true, 1, $false, "18+76=", 2) | Out-Null
$d.Content.Find.Execute("60+12=", $true, $false, $false, $false, $false, $true, 1, $false, "81-47=", 2) | Out-Null
$d.Content.Find.Execute("4+90=", $true, $false, $false, $false, $false, $true, 1, $false, "97-91=", 2) | Out-Null
$d.Content.Find.Execute("50+30=", $true, $false, $false, $false, $false, $true, 1, $false, "64-59=", 2) | Out-Null
$d.Content.Find.Execute("47+21=", $true, $false, $false, $false, $false, $true, 1, $false, "32+26=", 2) | Out-Null
$d.Content.Find.Execute("50-44=", $true, $false, $false, $false, $false, $true, 1, $false, "77-62=", 2) | Out-Null
$d.Content.Find.Execute("65+15=", $true, $false, $false, $false, $false, $true, 1, $false, "12-4=", 2) | Out-Null
$d.Content.Find.Execute("72-56=", $true, $false, $false, $false, $false, $true, 1, $false, "14+30=", 2) | Out-Null
$d.Content.Find.Execute("58-28=", $true, $false, $false, $false, $false, $true, 1, $false, "74-62=", 2) | Out-Null
$d.Content.Find.Execute("56-53=", $true, $false, $false, $false, $false, $true, 1, $false, "31+60=", 2) | Out-Null
$d.Content.Find.Execute("53-1=", $true, $false, $false, $false, $false, $true, 1, $false, "37+33=", 2) | Out-Null
$d.Content.Find.Execute("1+43=", $true, $false, $false, $false, $false, $true, 1, $false, "2+84=", 2) | Out-Null
$d.Content.Find.Execute("44-13=", $true, $false, $false, $false, $false, $true, 1, $false, "3+14=", 2) | Out-Null
$d.Content.Find.Execute("51-15=", $true, $false, $false, $false, $false, $true, 1, $false, "21+6=", 2) | Out-Null
$d.Content.Find.Execute("95-15=", $true, $false, $false, $false, $false, $true, 1, $false, "69+23=", 2) | Out-Null
$d.Content.Find.Execute("1+97=", $true, $false, $false, $false, $false, $true, 1, $false, "47+29=", 2) | Out-Null
$d.Content.Find.Execute("9+71=", $true, $false, $false, $false, $false, $true, 1, $false, "23-1=", 2) | Out-Null
$d.Content.Find.Execute("81-54=", $true, $false, $false, $false, $false, $true, 1, $false, "64-34=", 2) | Out-Null
$d.Content.Find.Execute("73-42=", $true, $false, $false, $false, $false, $true, 1, $false, "73-17=", 2) | Out-Null
$d.Content.Find.Execute("48-14=", $true, $false, $false, $false, $false, $true, 1, $false, "30+19=", 2) | Out-Null
$d.Content.Find.Execute("68-55=", $true, $false, $false, $false, $false, $true, 1, $false, "81+7=", 2) | Out-Null
$d.Content.Find.Execute("86-34=", $true, $false, $false, $false, $false, $true, 1, $false, "97-19=", 2) | Out-Null
$d.Content.Find.Execute("91-61=", $true, $false, $false, $false, $false, $true, 1, $false, "71+11=", 2) | Out-Null
$d.Content.Find.Execute("17+52=", $true, $false, $false, $false, $false, $true, 1, $false, "16+65=", 2) | Out-Null
$d.Content.Find.Execute("3+84=", $true, $false, $false, $false, $false, $true, 1, $false, "43+1=", 2) | Out-Null
$d.Content.Find.Execute("81-72=", $true, $false, $false, $false, $false, $true, 1, $false, "14+59=", 2) | Out-Null
$d.Content.Find.Execute("81-57=", $true, $false, $false, $false, $false, $true, 1, $false, "42+57=", 2) | Out-Null
$d.Content.Find.Execute("57-9=", $true, $false, $false, $false, $false, $true, 1, $false, "79+0=", 2) | Out-Null
$d.Content.Find.Execute("70-6=", $true, $false, $false, $false, $false, $true, 1, $false, "77+7=", 2) | Out-Null
$d.Content.Find.Execute("98-27=", $true, $false, $false, $false, $false, $true, 1, $false, "96-22=", 2) | Out-Null
$d.Content.Find.Execute("56+14=", $true, $false, $false, $false, $false, $true, 1, $false, "60+30=", 2) | Out-Null
$d.Content.Find.Execute("69-11=", $true, $false, $false, $false, $false, $true, 1, $false, "15-8=", 2) | Out-Null
$d.Content.Find.Execute("79-40=", $true, $false, $false, $false, $false, $true, 1, $false, "25-1=", 2) | Out-Null
$d.Content.Find.Execute("82-56=", $true, $false, $false, $false, $false, $true, 1, $false, "41+45=", 2) | Out-Null
$d.Content.Find.Execute("60+20=", $true, $false, $false, $false, $false, $true, 1, $false, "69-6=", 2) | Out-Null
$d.Content.Find.Execute("90-77=", $true, $false, $false, $false, $false, $true, 1, $false, "22-7=", 2) | Out-Null
$d.Content.Find.Execute("82-60=", $true, $false, $false, $false, $false, $true, 1, $false, "64-62=", 2) | Out-Null
$d.Content.Find.Execute("57-7=", $true, $false, $false, $false, $false, $true, 1, $false, "21+7=", 2) | Out-Null
$d.Content.Find.Execute("2+58=", $true, $false, $false, $false, $false, $true, 1, $false, "6+82=", 2) | Out-Null
$d.Content.Find.Execute("36-11=", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=", 2) | Out-Null
$d.Content.Find.Execute("62-53=", $true, $false, $false, $false, $false, $true, 1, $false, "59-9=", 2) | Out-Null
$d.Content.Find.Execute("88+10=", $true, $false, $false, $false, $false, $true, 1, $false, "41+52=", 2) | Out-Null
$d.Content.Find.Execute("59+20=", $true, $false, $false, $false, $false, $true, 1, $false, "90-46=", 2) | Out-Null
$d.Content.Find.Execute("99-93=", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=", 2) | Out-Null
$d.Content.Find.Execute("40-34=", $true, $false, $false, $false, $false, $true, 1, $false, "63-21=", 2) | Out-Null
$d.Content.Find.Execute("43-38=", $true, $false, $false, $false, $false, $true, 1, $false, "57-40=", 2) | Out-Null
$d.Content.Find.Execute("56-1=", $true, $false, $false, $false, $false, $true, 1, $false, "86-32=", 2) | Out-Null
$d.Content.Find.Execute("88-53=", $true, $false, $false, $false, $false, $true, 1, $false, "97-3=", 2) | Out-Null
$d.Content.Find.Execute("7+69=", $true, $false, $false, $false, $false, $true, 1, $false, "58-23=", 2) | Out-Null
$d.Content.Find.Execute("0+77=", $true, $false, $false, $false, $false, $true, 1, $false, "18+80=", 2) | Out-Null
$d.Content.Find.Execute("80-25=", $true, $false, $false, $false, $false, $true, 1, $false, "93-1=", 2) | Out-Null
$d.Content.Find.Execute("52-30=", $true, $false, $false, $false, $false, $true, 1, $false, "88-61=", 2) | Out-Null
$d.Content.Find.Execute("29+61=", $true, $false, $false, $false, $false, $true, 1, $false, "13+58=", 2) | Out-Null
$d.Content.Find.Execute("90-56=", $true, $false, $false, $false, $false, $true, 1, $false, "66+0=", 2) | Out-Null
$d.Content.Find.Execute("63-31=", $true, $false, $false, $false, $false, $true, 1, $false, "65+33=", 2) | Out-Null
$d.Content.Find.Execute("88-32=", $true, $false, $false, $false, $false, $true, 1, $false, "0+31=", 2) | Out-Null
$d.Content.Find.Execute("23+35=", $true, $false, $false, $false, $false, $true, 1, $false, "78-27=", 2) | Out-Null
$d.Content.Find.Execute("98-7=", $true, $false, $false, $false, $false, $true, 1, $false, "64-31=", 2) | Out-Null
$d.Content.Find.Execute("96-4=", $true, $false, $false, $false, $false, $true, 1, $false, "50-28=", 2) | Out-Null
$d.Content.Find.Execute("48-1=", $true, $false, $false, $false, $false, $true, 1, $false, "61-14=", 2) | Out-Null
$d.Content.Find.Execute("67-18=", $true, $false, $false, $false, $false, $true, 1, $false, "43-16=", 2) | Out-Null
$d.Content.Find.Execute("2+68=", $true, $false, $false, $false, $false, $true, 1, $false, "8+58=", 2) | Out-Null
$d.Content.Find.Execute("72-72=", $true, $false, $false, $false, $false, $true, 1, $false, "75-8=", 2) | Out-Null
$d.Content.Find.Execute("29+69=", $true, $false, $false, $false, $false, $true, 1, $false, "13+4=", 2) | Out-Null
$d.Content.Find.Execute("13+46=", $true, $false, $false, $false, $false, $true, 1, $false, "29+0=", 2) | Out-Null
$d.Content.Find.Execute("87-37=", $true, $false, $false, $false, $false, $true, 1, $false, "67-41=", 2) | Out-Null
$d.Content.Find.Execute("32+28=", $true, $false, $false, $false, $false, $true, 1, $false, "49+18=", 2) | Out-Null
$d.Content.Find.Execute("22-6=", $true, $false, $false, $false, $false, $true, 1, $false, "3+56=", 2) | Out-Null
$d.Content.Find.Execute("22+16=", $true, $false, $false, $false, $false, $true, 1, $false, "21+69=", 2) | Out-Null
$d.Content.Find.Execute("64+22=", $true, $false, $false, $false, $false, $true, 1, $false, "89-53=", 2) | Out-Null
$d.Content.Find.Execute("77-74=", $true, $false, $false, $false, $false, $true, 1, $false, "29-7=", 2) | Out-Null
$d.Content.Find.Execute("21-5=", $true, $false, $false, $false, $false, $true, 1, $false, "33+64=", 2) | Out-Null
$d.Content.Find.Execute("33-24=", $true, $false, $false, $false, $false, $true, 1, $false, "99+0=", 2) | Out-Null
$d.Content.Find.Execute("75-32=", $true, $false, $false, $false, $false, $true, 1, $false, "54+20=", 2) | Out-Null
$d.Content.Find.Execute("38-31=", $true, $false, $false, $false, $false, $true, 1, $false, "73-67=", 2) | Out-Null
$d.Content.Find.Execute("9+4=", $true, $false, $false, $false, $false, $true, 1, $false, "45-31=", 2) | Out-Null
$d.Content.Find.Execute("18+7=", $true, $false, $false, $false, $false, $true, 1, $false, "59-49=", 2) | Out-Null
$d.Content.Find.Execute("31+13=", $true, $false, $false, $false, $false, $true, 1, $false, "71-20=", 2) | Out-Null
$d.Content.Find.Execute("71-39=", $true, $false, $false, $false, $false, $true, 1, $false, "21+44=", 2) | Out-Null
